$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.889.20"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.45"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7535"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.41"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3038"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.37"
$ws.Range("E9").Value = "  -6.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06832"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07973"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.906.09"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7463"
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.199"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.17"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.895.92"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.91"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.940"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.19"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007718"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.937"
$ws.Range("E23").Value = "  +4.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.223"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.67"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1301"
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.016"
$ws.Range("E28").Value = "  -4.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.411"
$ws.Range("E29").Value = "  +4.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.518"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.275"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.019"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05351"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.249"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7246"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.718"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01913"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.791"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.174"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4400"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.13"
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.908"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8244"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.05"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.552"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.783"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.062.95"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.22"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05969"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.464"
$ws.Range("E51").Value = "  -0.31%  "
